$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.971.74"
$ws.Range("E2").Value = "  +2.96%  "

$ws.Range("D3").Value = "3.798.73"
$ws.Range("E3").Value = "  +0.96%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'704.89"
$ws.Range("E5").Value = "  +11.99%  "

$ws.Range("D6").Value = "'173.68"
$ws.Range("E6").Value = "  +5.04%  "

$ws.Range("D7").Value = "3.796.61"
$ws.Range("E7").Value = "  +0.99%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "'0.526"
$ws.Range("E9").Value = "  +1.22%  "

$ws.Range("E10").Value = "  +3.13%  "

$ws.Range("D11").Value = "'7.37"
$ws.Range("E11").Value = "  +8.89%  "

$ws.Range("D12").Value = "'0.463"
$ws.Range("E12").Value = "  +1.51%  "

$ws.Range("E13").Value = "  +7.70%  "

$ws.Range("D14").Value = "'36.41"
$ws.Range("E14").Value = "  +4.57%  "

$ws.Range("D15").Value = "4.435.73"
$ws.Range("E15").Value = "  +1.00%  "

$ws.Range("D16").Value = "3.796.54"
$ws.Range("E16").Value = "  +0.95%  "

$ws.Range("D17").Value = "71.018.03"
$ws.Range("E17").Value = "  +3.06%  "

$ws.Range("D18").Value = "'18.03"
$ws.Range("E18").Value = "  +2.25%  "

$ws.Range("D19").Value = "'7.27"
$ws.Range("E19").Value = "  +3.73%  "

$ws.Range("E20").Value = "  +1.00%  "

$ws.Range("D21").Value = "'11.26"
$ws.Range("E21").Value = "  +19.01%  "

$ws.Range("D22").Value = "'484.29"
$ws.Range("E22").Value = "  +4.99%  "

$ws.Range("D23").Value = "'0.721"
$ws.Range("E23").Value = "  +2.62%  "

$ws.Range("D24").Value = "'83.92"
$ws.Range("E24").Value = "  +2.38%  "

$ws.Range("D25").Value = "'0.0000146"
$ws.Range("E25").Value = "  +1.40%  "

$ws.Range("D26").Value = "'12.53"
$ws.Range("E26").Value = "  +3.40%  "

$ws.Range("D27").Value = "'10.62"
$ws.Range("E27").Value = "  +5.24%  "

$ws.Range("D28").Value = "'2.18"
$ws.Range("E28").Value = "  +2.96%  "

$ws.Range("D29").Value = "3.946.88"
$ws.Range("E29").Value = "  +0.90%  "

$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("D31").Value = "'3.13"
$ws.Range("E31").Value = "  +17.51%  "

$ws.Range("D32").Value = "'2.31"
$ws.Range("E32").Value = "  +1.84%  "

$ws.Range("D33").Value = "'7.61"
$ws.Range("E33").Value = "  +8.00%  "

$ws.Range("D34").Value = "'29.65"
$ws.Range("E34").Value = "  +4.77%  "

$ws.Range("E35").Value = "  +2.80%  "

$ws.Range("D36").Value = "'9.29"
$ws.Range("E36").Value = "  +4.77%  "

$ws.Range("E37").Value = "  +0.12%  "

$ws.Range("D38").Value = "3.746.08"
$ws.Range("E38").Value = "  +0.89%  "

$ws.Range("E39").Value = "  +3.30%  "

$ws.Range("E40").Value = "  +6.78%  "

$ws.Range("E41").Value = "  +4.29%  "

$ws.Range("D42").Value = "'2.25"
$ws.Range("E42").Value = "  +15.06%  "

$ws.Range("D43").Value = "'0.000332"
$ws.Range("E43").Value = "  +25.89%  "

$ws.Range("D44").Value = "'0.971"
$ws.Range("E44").Value = "  +1.27%  "

$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D47").Value = "'45.74"
$ws.Range("E47").Value = "  +7.35%  "

$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'49.49"
$ws.Range("E48").Value = "  +5.20%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'160.88"
$ws.Range("E49").Value = "  +2.54%  "

$ws.Range("E50").Value = "  -0.26%  "

$ws.Range("E51").Value = "  +2.64%  "
